# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{C="1017162179"; D="SANDRA JOHANA ACEVEDO VANEGAS"; E="1912"; F=33125; G=828116},
    @{C="1017162179"; D="SANDRA JOHANA ACEVEDO VANEGAS"; E="1911"; F=33125; G=828116},
    @{C="71314293";   D="JOHN FABER MARTINEZ ALZATE";    E="1912"; F=33125; G=828116},
    @{C="71314293";   D="JOHN FABER MARTINEZ ALZATE";    E="1911"; F=33125; G=828116},
    @{C="43186853";   D="BIBIANA ECHEVERRI RAMIREZ";     E="1912"; F=33125; G=828116},
    @{C="43186853";   D="BIBIANA ECHEVERRI RAMIREZ";     E="1911"; F=33125; G=828116},
    @{C="1128433590"; D="NATALI GIRALDO VALENCIA";       E="1912"; F=33125; G=828116},
    @{C="1128433590"; D="NATALI GIRALDO VALENCIA";       E="1911"; F=33125; G=828116},
    @{C="43108510";   D="TATIANA LICED ARDILA AYA";      E="1912"; F=33125; G=828116},
    @{C="43108510";   D="TATIANA LICED ARDILA AYA";      E="1911"; F=33125; G=828116},
    @{C="43622306";   D="SANDRA MILENA ACEVEDO MARIN";   E="1912"; F=40000; G=1000000},
    @{C="43622306";   D="SANDRA MILENA ACEVEDO MARIN";   E="1911"; F=40000; G=1000000},
    @{C="43163272";   D="CLAUDIA ANDREA VAHOS RODRIGUEZ";E="1912"; F=40000; G=1000000},
    @{C="43163272";   D="CLAUDIA ANDREA VAHOS RODRIGUEZ";E="1911"; F=40000; G=1000000}
)

$row = 16
foreach ($item in $data) {
    $ws.Cells.Item($row, 3).Value = $item.C
    $ws.Cells.Item($row, 4).Value = $item.D
    $ws.Cells.Item($row, 5).Value = $item.E
    $ws.Cells.Item($row, 6).Value = $item.F
    $ws.Cells.Item($row, 7).Value = $item.G
    $row++
}
